# Add a "Save" column (H) to the s_vals sheet, mirroring the header/style of
# the existing "sum" column (G) and filling in the per-row save flag values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1: copy formatting (bold font, border, centered/top alignment)
# from G1 so the new column matches the rest of the header row, then set the
# text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Per-row "Save" values for rows 2-53 (one entry per data row, in order).
$saveValues = @(
    0,0,0,0,0,0,0,1,0,0,
    1,0,0,0,0,0,0,0,0,0,
    0,0,0,0,0,0,0,0,0,0,
    0,1,0,1,0,0,0,0,0,0,
    0,0,1,0,1,0,0,0,0,0,
    0,0
)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
